$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 95, shifting the existing data (old rows 95..163)
# down to become rows 97..165.
$ws.Rows("95:96").Insert()

# Make sure the date column (D) keeps the date number format used elsewhere
# in the column (row 94, just above the inserted rows).
$dateFormat = $ws.Cells.Item(94, 4).NumberFormat
$ws.Cells.Item(95, 4).NumberFormat = $dateFormat
$ws.Cells.Item(96, 4).NumberFormat = $dateFormat

# --- Row 95 ---
$ws.Cells.Item(95, 1).Value = 11
$ws.Cells.Item(95, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(95, 3).Value = "Bíobío"
$ws.Cells.Item(95, 4).Value = 44484
$ws.Cells.Item(95, 5).Value = 8
$ws.Cells.Item(95, 6).Value = "Fruta"
$ws.Cells.Item(95, 7).Value = 100101
$ws.Cells.Item(95, 8).Value = "Berries"
$ws.Cells.Item(95, 9).Value = 100112025
$ws.Cells.Item(95, 10).Value = "Frutilla"
$ws.Cells.Item(95, 11).Value = "Sin especificar"
$ws.Cells.Item(95, 12).Value = "Especial"
$ws.Cells.Item(95, 13).Value = 100
$ws.Cells.Item(95, 14).Value = 12000
$ws.Cells.Item(95, 15).Value = 12000
$ws.Cells.Item(95, 16).Value = 12000
$ws.Cells.Item(95, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(95, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(95, 19).Value = 1714
$ws.Cells.Item(95, 20).Value = 7

# --- Row 96 ---
$ws.Cells.Item(96, 1).Value = 11
$ws.Cells.Item(96, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(96, 3).Value = "Bíobío"
$ws.Cells.Item(96, 4).Value = 44484
$ws.Cells.Item(96, 5).Value = 8
$ws.Cells.Item(96, 6).Value = "Fruta"
$ws.Cells.Item(96, 7).Value = 100101
$ws.Cells.Item(96, 8).Value = "Berries"
$ws.Cells.Item(96, 9).Value = 100112025
$ws.Cells.Item(96, 10).Value = "Frutilla"
$ws.Cells.Item(96, 11).Value = "Sin especificar"
$ws.Cells.Item(96, 12).Value = "Primera"
$ws.Cells.Item(96, 13).Value = 100
$ws.Cells.Item(96, 14).Value = 9000
$ws.Cells.Item(96, 15).Value = 9000
$ws.Cells.Item(96, 16).Value = 9000
$ws.Cells.Item(96, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(96, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(96, 19).Value = 1286
$ws.Cells.Item(96, 20).Value = 7
